$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 3228.1428
$ws.Range("I12").Value = 562
$ws.Range("J12").Value = 6783
$ws.Range("K12").Value = 562
$ws.Range("L12").Value = 6783
$ws.Range("M12").Value = -392
$ws.Range("N12").Value = -7123

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 922.381
$ws.Range("J17").Value = 922.35486
$ws.Range("L17").Value = 2767.06458
$ws.Range("N17").Value = -3103.06458

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1985.7778
$ws.Range("I32").Value = 1990
$ws.Range("J32").Value = 1971
$ws.Range("K32").Value = 1990
$ws.Range("L32").Value = 1971
$ws.Range("M32").Value = -1664
$ws.Range("N32").Value = -2623

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 133
$ws.Range("I33").Value = 113.27273
$ws.Range("K33").Value = 113.27273
$ws.Range("M33").Value = 115.72727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 14972.637
$ws.Range("I74").Value = 14972.637
$ws.Range("K74").Value = 14972.637
$ws.Range("M74").Value = -14036.637

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 14972.637
$ws.Range("I77").Value = 14972.637
$ws.Range("K77").Value = 74863.185
$ws.Range("M77").Value = -70183.185

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 919.64105
$ws.Range("I80").Value = 833.7
$ws.Range("J80").Value = 1010.1053
$ws.Range("K80").Value = 2501.1
$ws.Range("L80").Value = 3030.3159
$ws.Range("M80").Value = -1503.1
$ws.Range("N80").Value = -5026.3159

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 919.64105
$ws.Range("I83").Value = 833.7
$ws.Range("J83").Value = 1010.1053
$ws.Range("K83").Value = 7503.3
$ws.Range("L83").Value = 9090.947700000001
$ws.Range("M83").Value = -2511.3
$ws.Range("N83").Value = -19074.9477

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4833.778
$ws.Range("J88").Value = 3626
$ws.Range("L88").Value = 3626
$ws.Range("N88").Value = -4438

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4833.778
$ws.Range("J91").Value = 3626
$ws.Range("L91").Value = 3626
$ws.Range("N91").Value = -6434

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3682.0625
$ws.Range("J100").Value = 6100
$ws.Range("L100").Value = 6100
$ws.Range("N100").Value = -7182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1833.1305
$ws.Range("J112").Value = 1674.3529
$ws.Range("L112").Value = 5023.0587
$ws.Range("N112").Value = -7239.0587

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7631.72
$ws.Range("I137").Value = 1989.2142
$ws.Range("J137").Value = 14813.091
$ws.Range("K137").Value = 5967.642599999999
$ws.Range("L137").Value = 44439.273
$ws.Range("M137").Value = -3417.642599999999
$ws.Range("N137").Value = -49539.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8914.107
$ws.Range("I61").Value = 12960.2
$ws.Range("K61").Value = 12960.2
$ws.Range("M61").Value = -12748.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1224.1082
$ws.Range("I110").Value = 899.6923
$ws.Range("J110").Value = 1990.909
$ws.Range("K110").Value = 899.6923
$ws.Range("L110").Value = 1990.909
$ws.Range("M110").Value = 1145.3077
$ws.Range("N110").Value = -6080.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5111.967
$ws.Range("I132").Value = 3303.8235
$ws.Range("J132").Value = 7476.4614
$ws.Range("K132").Value = 9911.470499999999
$ws.Range("L132").Value = 22429.3842
$ws.Range("M132").Value = -7381.470499999999
$ws.Range("N132").Value = -27489.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8914.107
$ws.Range("I136").Value = 12960.2
$ws.Range("K136").Value = 38880.60000000001
$ws.Range("M136").Value = -36330.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25163
$ws.Range("I82").Value = 12097
$ws.Range("K82").Value = 12097
$ws.Range("M82").Value = -11714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 25163
$ws.Range("I85").Value = 12097
$ws.Range("K85").Value = 12097
$ws.Range("M85").Value = -10771

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1777.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6616.7334
$ws.Range("I134").Value = 7327.077
$ws.Range("K134").Value = 21981.231
$ws.Range("M134").Value = -19446.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4120.6
$ws.Range("I58").Value = 2316.1177
$ws.Range("J58").Value = 5824.8335
$ws.Range("K58").Value = 2316.1177
$ws.Range("L58").Value = 5824.8335
$ws.Range("M58").Value = -2113.1177
$ws.Range("N58").Value = -6230.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 33574.75
$ws.Range("I122").Value = 2096
$ws.Range("J122").Value = 65053.5
$ws.Range("K122").Value = 6288
$ws.Range("L122").Value = 195160.5
$ws.Range("M122").Value = -3838
$ws.Range("N122").Value = -200060.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3466.0881
$ws.Range("I134").Value = 3345.2334
$ws.Range("K134").Value = 10035.7002
$ws.Range("M134").Value = -7500.700199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4120.6
$ws.Range("I136").Value = 2316.1177
$ws.Range("J136").Value = 5824.8335
$ws.Range("K136").Value = 6948.353099999999
$ws.Range("L136").Value = 17474.5005
$ws.Range("M136").Value = -4398.353099999999
$ws.Range("N136").Value = -22574.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 48571.43
$ws.Range("J37").Value = 48571.43
$ws.Range("L37").Value = 145714.29
$ws.Range("N37").Value = -145938.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 812.36365
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 812.36365
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2437.09095
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6777.09095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 11227.417
$ws.Range("I124").Value = 8532.714
$ws.Range("K124").Value = 25598.142
$ws.Range("M124").Value = -20688.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 436935.1
$ws.Range("I129").Value = 1001738.8
$ws.Range("J129").Value = 2470.6924
$ws.Range("K129").Value = 3005216.4
$ws.Range("L129").Value = 7412.0772
$ws.Range("M129").Value = -3000216.4
$ws.Range("N129").Value = -17412.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4085.8
$ws.Range("I140").Value = 4085.8
$ws.Range("K140").Value = 12257.4
$ws.Range("M140").Value = -7077.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 859.5
$ws.Range("I16").Value = 859.5
$ws.Range("K16").Value = 859.5
$ws.Range("M16").Value = -689.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5261.478
$ws.Range("J132").Value = 8306.111000000001
$ws.Range("L132").Value = 24918.333
$ws.Range("N132").Value = -29978.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13159.583
$ws.Range("J62").Value = 19990.6
$ws.Range("L62").Value = 19990.6
$ws.Range("N62").Value = -21238.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 13159.583
$ws.Range("J65").Value = 19990.6
$ws.Range("L65").Value = 99953
$ws.Range("N65").Value = -106193

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1702.24
$ws.Range("I107").Value = 1666.5
$ws.Range("J107").Value = 1747.7273
$ws.Range("K107").Value = 4999.5
$ws.Range("L107").Value = 5243.1819
$ws.Range("M107").Value = -3079.5
$ws.Range("N107").Value = -9083.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52554.09
$ws.Range("I122").Value = 1175.7222
$ws.Range("K122").Value = 3527.1666
$ws.Range("M122").Value = -1077.1666
